# The diff just reorders <w:b/>/<w:i/> before <w:color/> inside <w:rPr> for a
# set of character styles (KeywordTok, ImportTok, CommentTok, ...), matching
# the element order required by wml.xsd. Re-assigning the Bold/Italic font
# property (even to its current value) causes the style's rPr to be
# re-serialized in schema order, which is exactly the fix described by the
# commit message/diff.

$d = $word.ActiveDocument

$boldStyles = @("KeywordTok", "ImportTok", "AnnotationTok", "CommentVarTok", "ControlFlowTok", "InformationTok", "WarningTok", "AlertTok", "ErrorTok")
$italicStyles = @("CommentTok", "DocumentationTok", "AnnotationTok", "CommentVarTok", "InformationTok", "WarningTok")

foreach ($name in $boldStyles) {
    $style = $d.Styles.Item($name)
    $font = $style.Font
    $font.Bold = $font.Bold
}

foreach ($name in $italicStyles) {
    $style = $d.Styles.Item($name)
    $font = $style.Font
    $font.Italic = $font.Italic
}
